$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6 (pushes old rows 6-8 down to 7-9)
$ws.Rows.Item(6).Insert()

# Fill in the new row's data: SecondGitHubRepo / ScondGithubRepo-TeamRead / pull
$ws.Range("A6").Value = "SecondGitHubRepo"
$ws.Range("B6").Value = "ScondGithubRepo-TeamRead"
$ws.Range("C6").Value = "pull"

# Update the active selection to C6 as reflected in the saved workbook view
$ws.Range("C6").Select()
